# Applies the commit: inserts a new lead row (A 61099-2024), refreshes the
# "Foerandrad" (C) date stamp to 45646 across all data rows, and appends two
# newly discovered notifications (A 60500-2024, A 60501-2024) at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 2; existing rows 2-34 shift down to 3-35 ---
$ws.Rows.Item(2).Insert()

# --- Fill the newly inserted row 2 with the new notification data ---
    $ws.Range("A2").Value2 = "A 61099-2024"
    $ws.Range("B2").Value2 = 45644
    $ws.Range("C2").Value2 = 45646
    $ws.Range("D2").Value2 = "OKÄNT"
    $ws.Range("E2").Value2 = "OKÄNT"
    $ws.Range("G2").Value2 = 0.4
    $ws.Range("H2").Value2 = 1
    $ws.Range("I2").Value2 = 0
    $ws.Range("J2").Value2 = 0
    $ws.Range("K2").Value2 = 0
    $ws.Range("L2").Value2 = 0
    $ws.Range("M2").Value2 = 0
    $ws.Range("N2").Value2 = 0
    $ws.Range("O2").Value2 = 0
    $ws.Range("P2").Value2 = 0
    $ws.Range("Q2").Value2 = 1
    $ws.Range("R2").Value2 = "Blåsippa"
    $ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/artfynd/A 61099-2024 artfynd.xlsx", "A 61099-2024")'
    $ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/kartor/A 61099-2024 karta.png", "A 61099-2024")'
    $ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomål/A 61099-2024 FSC-klagomål.docx", "A 61099-2024")'
    $ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomålsmail/A 61099-2024 FSC-klagomål mail.docx", "A 61099-2024")'
    $ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsyn/A 61099-2024 tillsynsbegäran.docx", "A 61099-2024")'
    $ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsynsmail/A 61099-2024 tillsynsbegäran mail.docx", "A 61099-2024")'

# Row 2 lost formatting when it was inserted blank: restore the date format
# on B/C and the wrap-text style on R to match the rest of the table.
$ws.Range("B2:C2").NumberFormat = "YYYY-MM-DD"
$ws.Range("R2").WrapText = $true

# --- Refresh the "Foerandrad" date (column C) to 45646 for every data row ---
for ($r = 2; $r -le 35; $r++) {
    $ws.Range("C$r").Value2 = 45646
}

# --- Append the two newly discovered notifications at the bottom (rows 36-37) ---

# Row 36: A 60500-2024
    $ws.Range("A36").Value2 = "A 60500-2024"
    $ws.Range("B36").Value2 = 45643
    $ws.Range("C36").Value2 = 45646
    $ws.Range("D36").Value2 = "OKÄNT"
    $ws.Range("E36").Value2 = "OKÄNT"
    $ws.Range("G36").Value2 = 0.8
    $ws.Range("H36").Value2 = 0
    $ws.Range("I36").Value2 = 0
    $ws.Range("J36").Value2 = 0
    $ws.Range("K36").Value2 = 0
    $ws.Range("L36").Value2 = 0
    $ws.Range("M36").Value2 = 0
    $ws.Range("N36").Value2 = 0
    $ws.Range("O36").Value2 = 0
    $ws.Range("P36").Value2 = 0
    $ws.Range("Q36").Value2 = 0
$ws.Range("B36:C36").NumberFormat = "YYYY-MM-DD"
$ws.Range("R36").WrapText = $true

# Row 37: A 60501-2024
    $ws.Range("A37").Value2 = "A 60501-2024"
    $ws.Range("B37").Value2 = 45643
    $ws.Range("C37").Value2 = 45646
    $ws.Range("D37").Value2 = "OKÄNT"
    $ws.Range("E37").Value2 = "OKÄNT"
    $ws.Range("G37").Value2 = 0.6
    $ws.Range("H37").Value2 = 0
    $ws.Range("I37").Value2 = 0
    $ws.Range("J37").Value2 = 0
    $ws.Range("K37").Value2 = 0
    $ws.Range("L37").Value2 = 0
    $ws.Range("M37").Value2 = 0
    $ws.Range("N37").Value2 = 0
    $ws.Range("O37").Value2 = 0
    $ws.Range("P37").Value2 = 0
    $ws.Range("Q37").Value2 = 0
$ws.Range("B37:C37").NumberFormat = "YYYY-MM-DD"
$ws.Range("R37").WrapText = $true
